$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update placeholder text to the new flattened variable names
$ws.Range("C2").Value = "Date: {{orderDate}}"
$ws.Range("C3").Value = "PO: #{{orderNumber}}"

$ws.Range("A7").Value = "{{vendorName}}"
$ws.Range("A8").Value = "{{vendorAddress}}"
$ws.Range("A9").Value = "{{vendorEmail}}"
$ws.Range("A10").Value = "{{vendorPhone}}"

$ws.Range("C7").Value = "{{customerName}}"
$ws.Range("C8").Value = "{{customerAddress}}"
$ws.Range("C9").Value = "{{customerEmail}}"
$ws.Range("C10").Value = "{{customerPhone}}"

# Update the active selection to match the saved view state
$ws.Range("A14").Select()
